$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new data rows, copying formatting from the last existing row pair ---
$ws.Range("C14:D15").Copy() | Out-Null
$ws.Range("C16:D17").PasteSpecial(-4122) | Out-Null

$ws.Range("C16").Value = "m.milasinovic@itsolivetti.it"
$ws.Range("D16").Value = 88
$ws.Range("C17").Value = "g.biancoli@itsolivetti.it"
$ws.Range("D17").Value = 90

# --- Selection like in target ---
$ws.Range("K6").Select() | Out-Null
